# This script rebuilds the workbook described by the diff:
#  - Sheet1 is repurposed to hold a PCB-congener-only summary (no ESTUARY
#    column): PARAMETRE_LIBELLE / First 5 years / Last 5 years (8 rows incl. header)
#  - A new Sheet2 is appended after Sheet1, holding the original
#    per-estuary breakdown but with the "median_1"/"median_2" headers
#    renamed to "First 5 years"/"Last 5 years" and new values.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# ---------------------------------------------------------------------
# Add Sheet2 right after Sheet1, *before* touching Sheet1's own data so
# Sheet1 stays the active sheet (tabSelected) exactly as before.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# Sheet2: per-estuary data (old content's layout), renamed headers +
# refreshed "First 5 years" / "Last 5 years" values.
# ---------------------------------------------------------------------
$sheet2Data = @(
    @("ESTUARY", "PARAMETRE_LIBELLE", "First 5 years", "Last 5 years"),
    @("Gironde", "CB 101", 157.119, 77.557),
    @("Loire",   "CB 101", 245.7235, 58.76000000000001),
    @("Seine",   "CB 101", 1805.3425, 569.072),
    @("Gironde", "CB 118", 58.799, 37.815),
    @("Loire",   "CB 118", 158.1035, 38.126),
    @("Seine",   "CB 118", 1533.4485, 373.984),
    @("Gironde", "CB 138", 527.2635, 128.1725),
    @("Loire",   "CB 138", 839.3395, 122.254),
    @("Seine",   "CB 138", 2982.968, 875.9014999999999),
    @("Gironde", "CB 153", 1122.002, 603.847),
    @("Loire",   "CB 153", 1126.1245, 310.6785),
    @("Seine",   "CB 153", 3999.73, 1530.8605),
    @("Gironde", "CB 180", 202.905, 67.52250000000001),
    @("Loire",   "CB 180", 107.019, 22.0375),
    @("Seine",   "CB 180", 369.52, 97.8295),
    @("Gironde", "CB 28", 12.531, 1.995),
    @("Loire",   "CB 28", 6.466, 1.6035),
    @("Seine",   "CB 28", 62.947, 12.094),
    @("Gironde", "CB 52", 51.0875, 19.163),
    @("Loire",   "CB 52", 47.8515, 11.8725),
    @("Seine",   "CB 52", 542.3315, 136.279)
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $r = $i + 1
    $row = $sheet2Data[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
}

$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("A1:D1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Sheet1: clear the old per-estuary content and replace with the
# congener-only summary (no ESTUARY column).
# ---------------------------------------------------------------------
$ws1.Cells.Clear()

$sheet1Data = @(
    @("PARAMETRE_LIBELLE", "First 5 years", "Last 5 years"),
    @("CB 101", 272.0365, 70.39949999999999),
    @("CB 118", 172.653, 40.273),
    @("CB 138", 892.4010000000001, 147.5985),
    @("CB 153", 1391.7845, 473.5635),
    @("CB 180", 202.905, 46.172),
    @("CB 28", 13.495, 1.9855),
    @("CB 52", 65.95699999999999, 14.42)
)

for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $r = $i + 1
    $row = $sheet1Data[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
}

$ws1.Range("A1:C1").Font.Bold = $true
$ws1.Range("A1:C1").HorizontalAlignment = -4108

$ws1.Activate()
